$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.34085166666667
$ws.Range("H2").Value = 43.022555
$ws.Range("I2").Value = 0.5393411052175457
$ws.Range("J2").Value = 0.5393411052175457
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.777551666666666
$ws.Range("N2").Value = 26.332655
$ws.Range("O2").Value = 0.07454818073713242
$ws.Range("P2").Value = 0.07454818073713242
$ws.Range("Q2").Value = 125.8775664481694
$ws.Range("R2").Value = 1132.898098033525
$ws.Range("S2").Value = 0.04020689819072235
$ws.Range("T2").Value = 0.04020689819072235
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.34085166666667
$ws.Range("H3").Value = 43.022555
$ws.Range("I3").Value = 0.5393411052175457
$ws.Range("J3").Value = 0.5393411052175457
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 36.14140700000001
$ws.Range("N3").Value = 108.424221
$ws.Range("O3").Value = 0.306950758417288
$ws.Range("P3").Value = 0.306950758417288
$ws.Range("Q3").Value = 518.2985568116284
$ws.Range("R3").Value = 4664.687011304655
$ws.Range("S3").Value = 0.165551161292144
$ws.Range("T3").Value = 0.165551161292144
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.34085166666667
$ws.Range("H4").Value = 43.022555
$ws.Range("I4").Value = 0.5393411052175457
$ws.Range("J4").Value = 0.5393411052175457
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 54.580447
$ws.Range("N4").Value = 163.741341
$ws.Range("O4").Value = 0.4635544377507104
$ws.Range("P4").Value = 0.4635544377507104
$ws.Range("Q4").Value = 782.7300943273616
$ws.Range("R4").Value = 7044.570848946255
$ws.Range("S4").Value = 0.2500139627849661
$ws.Range("T4").Value = 0.2500139627849661
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 14.34085166666667
$ws.Range("H5").Value = 43.022555
$ws.Range("I5").Value = 0.5393411052175457
$ws.Range("J5").Value = 0.5393411052175457
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 18.24393266666667
$ws.Range("N5").Value = 54.731798
$ws.Range("O5").Value = 0.1549466230948692
$ws.Range("P5").Value = 0.1549466230948692
$ws.Range("Q5").Value = 261.6335321893212
$ws.Range("R5").Value = 2354.70178970389
$ws.Range("S5").Value = 0.08356908294971327
$ws.Range("T5").Value = 0.08356908294971326
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.825018
$ws.Range("H6").Value = 20.475054
$ws.Range("I6").Value = 0.2566802053887532
$ws.Range("J6").Value = 0.2566802053887532
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.777551666666666
$ws.Range("N6").Value = 26.332655
$ws.Range("O6").Value = 0.07454818073713242
$ws.Range("P6").Value = 0.07454818073713242
$ws.Range("Q6").Value = 59.90694812093
$ws.Range("R6").Value = 539.16253308837
$ws.Range("S6").Value = 0.01913504234296505
$ws.Range("T6").Value = 0.01913504234296505
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.825018
$ws.Range("H7").Value = 20.475054
$ws.Range("I7").Value = 0.2566802053887532
$ws.Range("J7").Value = 0.2566802053887532
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 36.14140700000001
$ws.Range("N7").Value = 108.424221
$ws.Range("O7").Value = 0.306950758417288
$ws.Range("P7").Value = 0.306950758417288
$ws.Range("Q7").Value = 246.665753320326
$ws.Range("R7").Value = 2219.991779882934
$ws.Range("S7").Value = 0.07878818371478306
$ws.Range("T7").Value = 0.07878818371478305
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.825018
$ws.Range("H8").Value = 20.475054
$ws.Range("I8").Value = 0.2566802053887532
$ws.Range("J8").Value = 0.2566802053887532
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.580447
$ws.Range("N8").Value = 163.741341
$ws.Range("O8").Value = 0.4635544377507104
$ws.Range("P8").Value = 0.4635544377507104
$ws.Range("Q8").Value = 372.512533223046
$ws.Range("R8").Value = 3352.612799007414
$ws.Range("S8").Value = 0.1189852482907204
$ws.Range("T8").Value = 0.1189852482907204
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.825018
$ws.Range("H9").Value = 20.475054
$ws.Range("I9").Value = 0.2566802053887532
$ws.Range("J9").Value = 0.2566802053887532
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.24393266666667
$ws.Range("N9").Value = 54.731798
$ws.Range("O9").Value = 0.1549466230948692
$ws.Range("P9").Value = 0.1549466230948692
$ws.Range("Q9").Value = 124.515168840788
$ws.Range("R9").Value = 1120.636519567092
$ws.Range("S9").Value = 0.03977173104028477
$ws.Range("T9").Value = 0.03977173104028477
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.9237333333333333
$ws.Range("H10").Value = 2.7712
$ws.Range("I10").Value = 0.03474043024127374
$ws.Range("J10").Value = 0.03474043024127375
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.777551666666666
$ws.Range("N10").Value = 26.332655
$ws.Range("O10").Value = 0.07454818073713242
$ws.Range("P10").Value = 0.07454818073713242
$ws.Range("Q10").Value = 8.108117059555555
$ws.Range("R10").Value = 72.97305353599999
$ws.Range("S10").Value = 0.002589835872512216
$ws.Range("T10").Value = 0.002589835872512216
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.9237333333333333
$ws.Range("H11").Value = 2.7712
$ws.Range("I11").Value = 0.03474043024127374
$ws.Range("J11").Value = 0.03474043024127375
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 36.14140700000001
$ws.Range("N11").Value = 108.424221
$ws.Range("O11").Value = 0.306950758417288
$ws.Range("P11").Value = 0.306950758417288
$ws.Range("Q11").Value = 33.38502235946667
$ws.Range("R11").Value = 300.4652012352
$ws.Range("S11").Value = 0.01066360141030186
$ws.Range("T11").Value = 0.01066360141030186
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.9237333333333333
$ws.Range("H12").Value = 2.7712
$ws.Range("I12").Value = 0.03474043024127374
$ws.Range("J12").Value = 0.03474043024127375
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.580447
$ws.Range("N12").Value = 163.741341
$ws.Range("O12").Value = 0.4635544377507104
$ws.Range("P12").Value = 0.4635544377507104
$ws.Range("Q12").Value = 50.41777824213333
$ws.Range("R12").Value = 453.7600041792
$ws.Range("S12").Value = 0.01610408060771142
$ws.Range("T12").Value = 0.01610408060771143
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.9237333333333333
$ws.Range("H13").Value = 2.7712
$ws.Range("I13").Value = 0.03474043024127374
$ws.Range("J13").Value = 0.03474043024127375
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 18.24393266666667
$ws.Range("N13").Value = 54.731798
$ws.Range("O13").Value = 0.1549466230948692
$ws.Range("P13").Value = 0.1549466230948692
$ws.Range("Q13").Value = 16.85252873528889
$ws.Range("R13").Value = 151.6727586176
$ws.Range("S13").Value = 0.005382912350748239
$ws.Range("T13").Value = 0.005382912350748239
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.499973666666667
$ws.Range("H14").Value = 13.499921
$ws.Range("I14").Value = 0.1692382591524273
$ws.Range("J14").Value = 0.1692382591524273
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.777551666666666
$ws.Range("N14").Value = 26.332655
$ws.Range("O14").Value = 0.07454818073713242
$ws.Range("P14").Value = 0.07454818073713242
$ws.Range("Q14").Value = 39.49875135780611
$ws.Range("R14").Value = 355.488762220255
$ws.Range("S14").Value = 0.0126164043309328
$ws.Range("T14").Value = 0.0126164043309328
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.499973666666667
$ws.Range("H15").Value = 13.499921
$ws.Range("I15").Value = 0.1692382591524273
$ws.Range("J15").Value = 0.1692382591524273
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 36.14140700000001
$ws.Range("N15").Value = 108.424221
$ws.Range("O15").Value = 0.306950758417288
$ws.Range("P15").Value = 0.306950758417288
$ws.Range("Q15").Value = 162.6353797762824
$ws.Range("R15").Value = 1463.718417986541
$ws.Range("S15").Value = 0.05194781200005909
$ws.Range("T15").Value = 0.05194781200005909
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.499973666666667
$ws.Range("H16").Value = 13.499921
$ws.Range("I16").Value = 0.1692382591524273
$ws.Range("J16").Value = 0.1692382591524273
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 54.580447
$ws.Range("N16").Value = 163.741341
$ws.Range("O16").Value = 0.4635544377507104
$ws.Range("P16").Value = 0.4635544377507104
$ws.Range("Q16").Value = 245.6105742148956
$ws.Range("R16").Value = 2210.495167934061
$ws.Range("S16").Value = 0.07845114606731243
$ws.Range("T16").Value = 0.07845114606731245
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.499973666666667
$ws.Range("H17").Value = 13.499921
$ws.Range("I17").Value = 0.1692382591524273
$ws.Range("J17").Value = 0.1692382591524273
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 18.24393266666667
$ws.Range("N17").Value = 54.731798
$ws.Range("O17").Value = 0.1549466230948692
$ws.Range("P17").Value = 0.1549466230948692
$ws.Range("Q17").Value = 82.09721657643979
$ws.Range("R17").Value = 738.8749491879581
$ws.Range("S17").Value = 0.02622289675412295
$ws.Range("T17").Value = 0.02622289675412295
